$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 9
$ws.Range("C3").Value = 0.15
$ws.Range("C4").Value = 0.4
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 3000

$ws.Range("H11").Select()
